$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.091.36'
$ws.Range("E2").Value = '  +1.12%  '

$ws.Range("D3").Value = '3.502.53'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '602.56'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").Value = '175.45'
$ws.Range("E6").Value = '  +3.80%  '

$ws.Range("E7").Value = '  -1.28%  '

$ws.Range("D8").Value = '3.497.74'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").Value = '0.192'
$ws.Range("E10").Value = '  -0.26%  '

$ws.Range("D11").Value = '7.23'
$ws.Range("E11").Value = '  +9.00%  '

$ws.Range("D12").Value = '0.580'
$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("D13").Value = '46.16'
$ws.Range("E13").Value = '  -1.59%  '

$ws.Range("D14").Value = '0.0000274'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").Value = '4.070.58'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").Value = '8.27'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("D17").Value = '610.28'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("D18").Value = '3.506.38'
$ws.Range("E18").Value = '  +0.26%  '

$ws.Range("D19").Value = '70.216.75'
$ws.Range("E19").Value = '  +1.20%  '

$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("D21").Value = '17.32'
$ws.Range("E21").Value = '  +0.83%  '

$ws.Range("D22").Value = '0.875'
$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").Value = '8.96'
$ws.Range("E23").Value = '  -14.72%  '

$ws.Range("D24").Value = '97.60'
$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("D25").Value = '15.49'
$ws.Range("E25").Value = '  -1.09%  '

$ws.Range("E26").Value = '  -3.53%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").Value = '2.54'
$ws.Range("E28").Value = '  -1.85%  '

$ws.Range("D29").Value = '33.69'
$ws.Range("E29").Value = '  +2.07%  '

$ws.Range("D30").Value = '8.96'
$ws.Range("E30").Value = '  -2.73%  '

$ws.Range("B31").Value = 'Stacks'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D31").Value = '2.95'
$ws.Range("E31").Value = '  -3.74%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '7.98'
$ws.Range("E32").Value = '  -4.93%  '

$ws.Range("D33").Value = '634.58'
$ws.Range("E33").Value = '  +15.17%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").Value = '  -3.80%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '6.82'
$ws.Range("E35").Value = '  -0.42%  '

$ws.Range("D36").Value = '3.54'
$ws.Range("E36").Value = '  +2.44%  '

$ws.Range("D37").Value = '0.0989'
$ws.Range("E37").Value = '  -1.84%  '

$ws.Range("D38").Value = '10.71'
$ws.Range("E38").Value = '  -0.18%  '

$ws.Range("D39").Value = '0.0469'
$ws.Range("E39").Value = '  +4.55%  '

$ws.Range("D40").Value = '56.62'
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("E42").Value = '  +1.90%  '

$ws.Range("D43").Value = '3.351.38'
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").Value = '0.0₃0733'
$ws.Range("E44").Value = '  +5.66%  '

$ws.Range("E45").Value = '  -5.40%  '

$ws.Range("D46").Value = '32.09'
$ws.Range("E46").Value = '  -2.51%  '

$ws.Range("D47").Value = '2.87'
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("E48").Value = '  -2.81%  '

$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").Value = '133.94'
$ws.Range("E50").Value = '  -0.36%  '
